$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: result_spell IDs in column E (rows 4-13) were off by a digit
# (202xx instead of 201xx).
$ws.Range("E4").Value = 20101
$ws.Range("E5").Value = 20102
$ws.Range("E6").Value = 20103
$ws.Range("E7").Value = 20104
$ws.Range("E8").Value = 20105
$ws.Range("E9").Value = 20106
$ws.Range("E10").Value = 20107
$ws.Range("E11").Value = 20108
$ws.Range("E12").Value = 20109
$ws.Range("E13").Value = 20110

# Update the view state: move the active selection (and with it the
# frozen-pane top-left cell) to E21.
[void]$ws.Range("E21").Select()

# Page setup: paper size 9 (A4), portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
